$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trial 6 row (row 9): change trial number, add epoch/accuracy/notes
$ws.Range("A9").Value = 6.5
$ws.Range("G9").Value = 15
$ws.Range("H9").Value = "~66%"
$ws.Range("I9").Value = "Reached best pretty soon and started to bounce around"

# Add new trial 7 row (row 10)
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Simple_MLP([40*(2*context_size+1), 256, 128, 128, 71])"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = "Adam"
$ws.Range("E10").Value = 0.0002
$ws.Range("F10").Value = 256

# Match number format of the Initial Learning Rate column (scientific notation)
$ws.Range("E10").NumberFormat = "0.00E+00"

# Update selection to reflect new active cell
$ws.Range("H10").Select()
